$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before existing row 3; this shifts old rows 3,4,5 down to 4,5,6
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new data set
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = 44624
$ws.Cells.Item(3, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = 100112044
$ws.Cells.Item(3, 7).Value = "Perejil"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 120
$ws.Cells.Item(3, 11).Value = 650
$ws.Cells.Item(3, 12).Value = 700
$ws.Cells.Item(3, 13).Value = 675
$ws.Cells.Item(3, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(3, 15).Value = "Región del Maule"
$ws.Cells.Item(3, 16).Value = 675
$ws.Cells.Item(3, 17).Value = 1
$ws.Cells.Item(3, 18).Value = "Hortaliza"
